$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PREGRADO")
$ws.Columns("L:L").Delete()
try {
  $ws.AutoFilter.Range = $ws.Range("A5:K7")
  Write-Output "SET OK"
} catch {
  Write-Output "ERR $_"
}
Write-Output $ws.AutoFilter.Range.Address()
